$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at H (column 8) for the "Open" field, shifting
# secType/expiration/strike/right/conID one column to the right.
$ws.Columns.Item(8).Insert()

# Populate the new "Open" column header and values.
$ws.Cells.Item(1, 8).Value = "Open"
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(3, 8).Value = 1

# Move the active selection to match the new layout.
[void]$ws.Range("H4").Select()

# Reflect the updated window geometry recorded by Excel on save.
$win = $wb.Windows.Item(1)
$win.Left = 6140
$win.Top = 460
$win.Width = 22660
$win.Height = 10000
